$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for a new data row right after the current last row (36), pushing
# nothing else around (it's the new last row of the table).
$ws.Rows.Item(37).Insert(-4121)

# Populate the new row's cells. The order below matters: it reproduces the
# order in which the new shared strings were appended to sharedStrings.xml
# in the target workbook (the "I'm having so much fun! ♪" English line first,
# then the script filename, then the Russian translation, then the
# 1-byte-shifted/"encoded" Russian line).
$ws.Range("C37").Value2 = " I\'m having so much fun! ♪"
$ws.Range("A37").Value2 = "SCRIPT/G01P03A/us2306.ssb"
$ws.Range("B37").Value2 = 19
$ws.Range("D37").Value2 = " Мне так весело! ♪"
$ws.Range("E37").Value2 = " Íîå óàë âåòåìï! ♪"

# New row keeps the same row height (43.2pt) as the rest of the table's
# entries.
$ws.Rows.Item(37).RowHeight = 43.2

# Row 36 now becomes the first row of a new "file name" group, so it gets the
# thin-top-and-bottom-border style used by the other group-leading rows
# (e.g. row 33) instead of the plain "continuation" style it had before.
$groupHeaderFormat = $ws.Range("A33:E33")
$groupHeaderFormat.Copy()
$ws.Range("A36:E36").PasteSpecial(-4122)

# Match the workbook's new selection/cursor position.
$ws.Range("C39").Select()
